# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The sheet lists two workers' overdue-payment periods. Previously the rows
# were grouped by worker (all of worker A's periods, then all of worker B's
# periods). This edit re-sorts the 14 data rows (16-29) so the two workers'
# records are interleaved period-by-period in chronological order
# (2312, 2401, 2402, 2403, 2404, 2405, 2406), alternating between the two
# workers, while keeping each worker's own Valor Mora / Salario Basico
# values attached to the correct period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2312", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2312", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2401", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2401", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2402", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2402", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2403", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2403", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2404", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2404", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2405", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2405", 92800, 2320000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",     "2406", 18560, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA", "2406", 37120, 2320000)
)

# Column layout: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 2).Value = $vals[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $vals[4]   # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $vals[5]   # G: Salario Basico
}
